# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The "Periodo Mora" column (E16:E29) is reordered so the periods run in
# ascending order (1908 .. 2009) instead of descending (2009 .. 1908), and
# the "Valor Mora" (F column) amount that used to sit on the 2009 row now
# sits on the 1908 row (i.e. the two distinct amounts swap rows along with
# the reordering).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New ascending order for the "Periodo Mora" column, rows 16-29.
$periods = @("1908","1909","1910","1911","1912","2001","2002","2003","2004","2005","2006","2007","2008","2009")

for ($i = 0; $i -lt $periods.Length; $i++) {
    $row = 16 + $i
    $ws.Cells.Item($row, 5).Value = $periods[$i]
}

# The distinct "Valor Mora" amount (23187) moves from the last row (29) to
# the first row (16); every other row keeps the common amount (33125).
$ws.Cells.Item(16, 6).Value = 23187
$ws.Cells.Item(29, 6).Value = 33125
